# Apply updates described in the commit diff (gh-pages data refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 130
$ws.Range("F5").Value = 361
$ws.Range("C6").Value = "上海·奇卡波利国潮嘉年华"
$ws.Range("F6").Value = 757
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202402/QBcbo0Do1707295657878.jpeg"
$ws.Range("F7").Value = 203
$ws.Range("F8").Value = 1068
$ws.Range("F9").Value = 269
$ws.Range("F11").Value = 350
$ws.Range("F12").Value = 619
$ws.Range("F14").Value = 492
$ws.Range("F18").Value = 821
$ws.Range("F21").Value = 24
$ws.Range("F24").Value = 204
$ws.Range("F26").Value = 150
$ws.Range("F28").Value = 959
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 221
$ws.Range("F32").Value = 1020

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1040
$ws.Range("F5").Value = 1040
$ws.Range("F8").Value = 225
$ws.Range("F9").Value = 17
$ws.Range("F14").Value = 583
$ws.Range("F17").Value = 968
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 290
$ws.Range("F25").Value = 260
$ws.Range("F26").Value = 3665
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 18
$ws.Range("F31").Value = 22
$ws.Range("F33").Value = 103

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2399
$ws.Range("F6").Value = 985
$ws.Range("F7").Value = 2
$ws.Range("F9").Value = 1241

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2399
$ws.Range("F7").Value = 985
$ws.Range("F8").Value = 1241
$ws.Range("F11").Value = 130
$ws.Range("F12").Value = 361
$ws.Range("C13").Value = "上海·奇卡波利国潮嘉年华"
$ws.Range("F13").Value = 757
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202402/QBcbo0Do1707295657878.jpeg"
$ws.Range("F14").Value = 203
$ws.Range("F16").Value = 1068
$ws.Range("F17").Value = 269
$ws.Range("F18").Value = 350
$ws.Range("F19").Value = 619
$ws.Range("F20").Value = 1040
$ws.Range("F21").Value = 492
$ws.Range("F24").Value = 821
$ws.Range("F27").Value = 24
$ws.Range("F30").Value = 204
$ws.Range("F31").Value = 150
$ws.Range("F34").Value = 959
$ws.Range("F35").Value = 583
$ws.Range("F36").Value = 583
$ws.Range("F39").Value = 221
$ws.Range("F44").Value = 290
$ws.Range("F45").Value = 290
$ws.Range("F46").Value = 260
$ws.Range("F47").Value = 1020
